$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns that would otherwise be auto-converted to numbers to remain text,
# matching the workbook's existing text-cell convention for the Price column.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "59.237.32"
$ws.Range("E2").Value = "  +1.16%  "
$ws.Range("D3").Value = "2.994.01"
$ws.Range("E3").Value = "  +0.27%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "560.88"
$ws.Range("E5").Value = "  -0.40%  "
$ws.Range("D6").Value = "138.02"
$ws.Range("E6").Value = "  +2.36%  "
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").Value = "2.983.51"
$ws.Range("E9").Value = "  +0.12%  "
$ws.Range("E11").Value = "  +5.14%  "
$ws.Range("E12").Value = "  +1.79%  "
$ws.Range("E13").Value = "  +0.82%  "
$ws.Range("E14").Value = "  +0.79%  "
$ws.Range("E15").Value = "  +1.37%  "
$ws.Range("D16").Value = "3.488.78"
$ws.Range("E16").Value = "  +0.31%  "
$ws.Range("E17").Value = "  +6.05%  "
$ws.Range("D18").Value = "2.991.01"
$ws.Range("E18").Value = "  +0.46%  "
$ws.Range("D19").Value = "59.243.00"
$ws.Range("E19").Value = "  +1.53%  "
$ws.Range("D20").Value = "430.35"
$ws.Range("E20").Value = "  +1.17%  "
$ws.Range("E21").Value = "  +1.69%  "
$ws.Range("D22").Value = "0.721"
$ws.Range("E22").Value = "  +4.04%  "
$ws.Range("D23").Value = "13.52"
$ws.Range("E23").Value = "  +2.08%  "
$ws.Range("E24").Value = "  +0.98%  "
$ws.Range("D25").Value = "80.29"
$ws.Range("E25").Value = "  +0.15%  "
$ws.Range("E26").Value = "  -0.17%  "
$ws.Range("D27").Value = "2.21"
$ws.Range("E27").Value = "  +8.80%  "
$ws.Range("E28").Value = "  +0.27%  "
$ws.Range("E29").Value = "  +0.52%  "
$ws.Range("D30").Value = "7.83"
$ws.Range("E30").Value = "  +2.11%  "
$ws.Range("D31").Value = "25.70"
$ws.Range("E31").Value = "  +0.48%  "
$ws.Range("E32").Value = "  -1.07%  "
$ws.Range("E33").Value = "  +0.55%  "
$ws.Range("D34").Value = "0.996"
$ws.Range("E34").Value = "  +4.69%  "
$ws.Range("E35").Value = "  +3.80%  "
$ws.Range("D36").Value = "0.0₃0756"
$ws.Range("E36").Value = "  +7.66%  "
$ws.Range("D37").Value = "2.09"
$ws.Range("E37").Value = "  -2.45%  "
$ws.Range("D38").Value = "48.89"
$ws.Range("E38").Value = "  +0.19%  "
$ws.Range("E39").Value = "  -1.23%  "
$ws.Range("E40").Value = "  +4.50%  "
$ws.Range("D41").Value = "404.59"
$ws.Range("E41").Value = "  +5.74%  "
$ws.Range("D42").Value = "2.771.06"
$ws.Range("E42").Value = "  +1.52%  "
$ws.Range("E43").Value = "  -0.39%  "
$ws.Range("E44").Value = "  -1.96%  "
$ws.Range("E45").Value = "  +2.96%  "
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("D47").Value = "123.25"
$ws.Range("E47").Value = "  +0.00%  "
$ws.Range("D48").Value = "34.32"
$ws.Range("E48").Value = "  +19.16%  "
$ws.Range("E49").Value = "  -0.72%  "
$ws.Range("E50").Value = "  -0.77%  "
$ws.Range("D51").Value = "23.44"
$ws.Range("E51").Value = "  -1.38%  "
